$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I1").Value = "21/03/2023"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 497
$ws.Range("D2").Value = 512
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 707.2
$ws.Range("J2").Value = -27.60180995475113
$ws.Range("C3").Value = 51
$ws.Range("D3").Value = 53
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("I3").Value = 81
$ws.Range("J3").Value = -34.5679012345679
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 137
$ws.Range("D4").Value = 152
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3
$ws.Range("I4").Value = 120
$ws.Range("J4").Value = 26.66666666666666
$ws.Range("C5").Value = 310
$ws.Range("D5").Value = 319
$ws.Range("E5").Value = 0
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 8
$ws.Range("I5").Value = 504
$ws.Range("J5").Value = -36.70634920634921
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 228
$ws.Range("D6").Value = 262
$ws.Range("E6").Value = 29
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 2
$ws.Range("I6").Value = 364
$ws.Range("J6").Value = -28.02197802197802
$ws.Range("C7").Value = 80
$ws.Range("D7").Value = 82
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 110
$ws.Range("J7").Value = -25.45454545454545
$ws.Range("C8").Value = 89
$ws.Range("D8").Value = 113
$ws.Range("E8").Value = 21
$ws.Range("F8").Value = 3
$ws.Range("I8").Value = 119
$ws.Range("J8").Value = -5.042016806722693
$ws.Range("C9").Value = 161
$ws.Range("D9").Value = 161
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 441
$ws.Range("J9").Value = -63.49206349206349
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 30
$ws.Range("D10").Value = 30
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 2
$ws.Range("I10").Value = 50
$ws.Range("J10").Value = -40
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 22
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 1
$ws.Range("I12").Value = 44
$ws.Range("J12").Value = -50
